$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data cells we are about to write stay text-typed (matches the
# original inlineStr/text cells) instead of Excel auto-converting
# number-looking strings (e.g. "212.45") into numeric cells.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.892.71"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.617.29"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "212.45"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "18.28"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "1.841.89"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "1.606.71"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").Value = "4.13"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "25.894.93"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "61.37"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "0.0₃0735"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "190.98"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("D25").Value = "143.75"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").Value = "15.21"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "1.22"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("D33").Value = "3.08"
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").Value = "1.128.89"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "0.839"
$ws.Range("E37").Value = "  -4.55%  "
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0153"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "0.510"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").Value = "97.84"
$ws.Range("D42").Value = "1.752.98"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Value = "0.747"
$ws.Range("E43").Value = "  -4.79%  "
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").Value = "53.94"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "0.410"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.46"
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.68%  "
